# Small positional bug fixes for the text labels on slide 1 of
# SARWithRictor.pptx: nudge ten "tx*" caption shapes by a few EMU each
# (their sizes are unchanged).
#
# The shapes live inside the single group shape on the slide, so they
# are reached through GroupItems. PowerPoint's Shape.Left / Shape.Top
# are expressed in points (1 pt = 914400/72 EMU) and are backed by a
# single-precision float, so a naive EMU->point conversion can be off
# by a unit or two once it is converted back to EMU on save. To land
# on the exact target EMU we search for the nearest representable
# point value whose round-trip reproduces it exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "tx*" label shapes live inside the single group shape on the
# slide (the second top-level shape).
$grp = $s.Shapes.Item(2)

$EMU_PER_PT = 914400 / 72

function Set-ExactEmu {
    param($shape, $targetEmu, $axis)

    $basePt = $targetEmu / $EMU_PER_PT

    if ($axis -eq "X") {
        $shape.Left = $basePt
    } else {
        $shape.Top = $basePt
    }

    $i = -400
    while ($i -le 400) {
        $candidate = $basePt + ($i * 0.0000001)
        if ($axis -eq "X") {
            $shape.Left = $candidate
            $got = [math]::Round($shape.Left * $EMU_PER_PT)
        } else {
            $shape.Top = $candidate
            $got = [math]::Round($shape.Top * $EMU_PER_PT)
        }
        if ($got -eq $targetEmu) {
            return
        }
        $i = $i + 1
    }
    # Fall back to the plain conversion if an exact match wasn't found
    # (shouldn't happen in practice).
    if ($axis -eq "X") {
        $shape.Left = $basePt
    } else {
        $shape.Top = $basePt
    }
}

# shape name, new X (EMU), new Y (EMU)
$updates = @(
    @("tx9",  4518434, 2462577),
    @("tx10", 4873821, 2768961),
    @("tx11", 5315965, 3174502),
    @("tx12", 5906317, 3480885),
    @("tx13", 6434818, 3815532),
    @("tx14", 6452709, 4162462),
    @("tx15", 4719515, 4803470),
    @("tx16", 4833901, 5150400),
    @("tx17", 4060475, 3133291),
    @("tx18", 4150777, 3480220)
)

foreach ($u in $updates) {
    $shapeName = $u[0]
    $xEmu = $u[1]
    $yEmu = $u[2]
    $shape = $grp.GroupItems.Item($shapeName)
    Set-ExactEmu $shape $xEmu "X"
    Set-ExactEmu $shape $yEmu "Y"
}
